$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# --- Simple value edits -------------------------------------------------
# End Year: 2090 -> 2060
$ws.Range("B4").Value = 2060

# investment_initialization_years: 0 -> 3
$ws.Range("B18").Value = 3

# increase_demand: FALSE -> TRUE (also gets highlighted like the other
# "active" switches around it, same yellow fill used on B28/B29/B30)
$ws.Range("B31").Value = $true
$ws.Range("B31").Interior.Color = $ws.Range("B28").Interior.Color

# --- Column width ---------------------------------------------------------
$ws.Columns("B").ColumnWidth = 38.5546875

# --- Move the "Checking inputs" block down by 2 rows (and merge the lone
# label row into the first formula row) -----------------------------------
$ws.Range("A46:C51").ClearContents()

$ws.Range("A49").Value = "Checking inputs"
$ws.Range("B49").Formula = '=IF(AND(B26=TRUE,B24>0),"PRICES are fixed, no fuel trends are considered","ok")'

$ws.Range("B50").Formula = '=IF(AND(B20=TRUE,B19=FALSE),"DANGER!!!!!","ok")'
$ws.Range("C50").Value = "If the dummy capacity will be installed, it could be very different than expected "

$ws.Range("B51").Formula = '=IF(AND(B20=FALSE,B19=TRUE),"DANGER","ok")'
$ws.Range("C51").Value = "Testing different capacity than the one being installed can cause deviations in reality"

$ws.Range("B52").Formula = '=IF(AND(B27=TRUE,B26=TRUE),"DANGER","ok")'
$ws.Range("C52").Value = "Testing different capacity than the one being installed can cause deviations in reality"

$ws.Range("B53").Formula = '=IF(AND(B41<>"NOTSET",B40<>"NOTSET"),"Either NPV or IRR","ok")'
$ws.Range("C53").Value = "Either IRR or NPV must be None"

# Conditional formatting tracked the old B47:B51 block - move it along.
$fc = $ws.Range("B47:B51").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("B49:B53"))

# --- View state: scroll back to top, select C5 instead of B45 -------------
$ws.Range("A1").Select()
$ws.Range("C5").Select()
